# Natmi following Dr Hou advice
# Adds a new "FAPs" cluster category and expands the LR-pairs results
# table from 2 data rows to 6 data rows (rows 2-7), covering all
# sending/target cluster combinations among FAPs / sCs / ECs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs -> ECs -------------------------------------------------
$ws.Cells.Item(2, 1).Value  = "FAPs"
$ws.Cells.Item(2, 2).Value  = "Fgf16"
$ws.Cells.Item(2, 3).Value  = "Fgfr3"
$ws.Cells.Item(2, 4).Value  = "ECs"
$ws.Cells.Item(2, 5).Value  = 2
$ws.Cells.Item(2, 6).Value  = 0.6666666666666666
$ws.Cells.Item(2, 7).Value  = 0.09755033333333334
$ws.Cells.Item(2, 8).Value  = 0.292651
$ws.Cells.Item(2, 9).Value  = 0.1198375636346959
$ws.Cells.Item(2, 10).Value = 0.1198375636346959
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 5.243417666666667
$ws.Cells.Item(2, 14).Value = 15.730253
$ws.Cells.Item(2, 15).Value = 0.8253998362974575
$ws.Cells.Item(2, 16).Value = 0.8253998362974574
$ws.Cells.Item(2, 17).Value = 0.5114971411892223
$ws.Cells.Item(2, 18).Value = 4.603474270703
$ws.Cells.Item(2, 19).Value = 0.09891390540636416
$ws.Cells.Item(2, 20).Value = 0.09891390540636412

# --- Row 3: FAPs -> FAPs -------------------------------------------------
$ws.Cells.Item(3, 1).Value  = "FAPs"
$ws.Cells.Item(3, 2).Value  = "Fgf16"
$ws.Cells.Item(3, 3).Value  = "Fgfr3"
$ws.Cells.Item(3, 4).Value  = "FAPs"
$ws.Cells.Item(3, 5).Value  = 2
$ws.Cells.Item(3, 6).Value  = 0.6666666666666666
$ws.Cells.Item(3, 7).Value  = 0.09755033333333334
$ws.Cells.Item(3, 8).Value  = 0.292651
$ws.Cells.Item(3, 9).Value  = 0.1198375636346959
$ws.Cells.Item(3, 10).Value = 0.1198375636346959
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.6792986666666666
$ws.Cells.Item(3, 14).Value = 2.037896
$ws.Cells.Item(3, 15).Value = 0.1069327381315001
$ws.Cells.Item(3, 16).Value = 0.1069327381315001
$ws.Cells.Item(3, 17).Value = 0.06626581136622221
$ws.Cells.Item(3, 18).Value = 0.5963923022959999
$ws.Cells.Item(3, 19).Value = 0.01281455881046591
$ws.Cells.Item(3, 20).Value = 0.01281455881046591

# --- Row 4: FAPs -> sCs --------------------------------------------------
$ws.Cells.Item(4, 1).Value  = "FAPs"
$ws.Cells.Item(4, 2).Value  = "Fgf16"
$ws.Cells.Item(4, 3).Value  = "Fgfr3"
$ws.Cells.Item(4, 4).Value  = "sCs"
$ws.Cells.Item(4, 5).Value  = 2
$ws.Cells.Item(4, 6).Value  = 0.6666666666666666
$ws.Cells.Item(4, 7).Value  = 0.09755033333333334
$ws.Cells.Item(4, 8).Value  = 0.292651
$ws.Cells.Item(4, 9).Value  = 0.1198375636346959
$ws.Cells.Item(4, 10).Value = 0.1198375636346959
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4298626666666667
$ws.Cells.Item(4, 14).Value = 1.289588
$ws.Cells.Item(4, 15).Value = 0.06766742557104236
$ws.Cells.Item(4, 16).Value = 0.06766742557104236
$ws.Cells.Item(4, 17).Value = 0.04193324642088889
$ws.Cells.Item(4, 18).Value = 0.3773992177880001
$ws.Cells.Item(4, 19).Value = 0.008109099417865837
$ws.Cells.Item(4, 20).Value = 0.008109099417865836

# --- Row 5: sCs -> ECs ----------------------------------------------------
$ws.Cells.Item(5, 1).Value  = "sCs"
$ws.Cells.Item(5, 2).Value  = "Fgf16"
$ws.Cells.Item(5, 3).Value  = "Fgfr3"
$ws.Cells.Item(5, 4).Value  = "ECs"
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 6).Value  = 1
$ws.Cells.Item(5, 7).Value  = 0.716471
$ws.Cells.Item(5, 8).Value  = 2.149413
$ws.Cells.Item(5, 9).Value  = 0.8801624363653041
$ws.Cells.Item(5, 10).Value = 0.880162436365304
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.243417666666667
$ws.Cells.Item(5, 14).Value = 15.730253
$ws.Cells.Item(5, 15).Value = 0.8253998362974575
$ws.Cells.Item(5, 16).Value = 0.8253998362974574
$ws.Cells.Item(5, 17).Value = 3.756756699054333
$ws.Cells.Item(5, 18).Value = 33.810810291489
$ws.Cells.Item(5, 19).Value = 0.7264859308910934
$ws.Cells.Item(5, 20).Value = 0.7264859308910931

# --- Row 6: sCs -> FAPs ----------------------------------------------------
$ws.Cells.Item(6, 1).Value  = "sCs"
$ws.Cells.Item(6, 2).Value  = "Fgf16"
$ws.Cells.Item(6, 3).Value  = "Fgfr3"
$ws.Cells.Item(6, 4).Value  = "FAPs"
$ws.Cells.Item(6, 5).Value  = 3
$ws.Cells.Item(6, 6).Value  = 1
$ws.Cells.Item(6, 7).Value  = 0.716471
$ws.Cells.Item(6, 8).Value  = 2.149413
$ws.Cells.Item(6, 9).Value  = 0.8801624363653041
$ws.Cells.Item(6, 10).Value = 0.880162436365304
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.6792986666666666
$ws.Cells.Item(6, 14).Value = 2.037896
$ws.Cells.Item(6, 15).Value = 0.1069327381315001
$ws.Cells.Item(6, 16).Value = 0.1069327381315001
$ws.Cells.Item(6, 17).Value = 0.4866977950053333
$ws.Cells.Item(6, 18).Value = 4.380280155048
$ws.Cells.Item(6, 19).Value = 0.09411817932103415
$ws.Cells.Item(6, 20).Value = 0.09411817932103414

# --- Row 7: sCs -> sCs ------------------------------------------------------
$ws.Cells.Item(7, 1).Value  = "sCs"
$ws.Cells.Item(7, 2).Value  = "Fgf16"
$ws.Cells.Item(7, 3).Value  = "Fgfr3"
$ws.Cells.Item(7, 4).Value  = "sCs"
$ws.Cells.Item(7, 5).Value  = 3
$ws.Cells.Item(7, 6).Value  = 1
$ws.Cells.Item(7, 7).Value  = 0.716471
$ws.Cells.Item(7, 8).Value  = 2.149413
$ws.Cells.Item(7, 9).Value  = 0.8801624363653041
$ws.Cells.Item(7, 10).Value = 0.880162436365304
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.4298626666666667
$ws.Cells.Item(7, 14).Value = 1.289588
$ws.Cells.Item(7, 15).Value = 0.06766742557104236
$ws.Cells.Item(7, 16).Value = 0.06766742557104236
$ws.Cells.Item(7, 17).Value = 0.3079841346493334
$ws.Cells.Item(7, 18).Value = 2.771857211844
$ws.Cells.Item(7, 19).Value = 0.05955832615317652
$ws.Cells.Item(7, 20).Value = 0.05955832615317651
